$d = $word.ActiveDocument
$d.Content.Find.Execute("Version 1.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version 2.", 2)
